$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking Accuracy/F1_Score columns keep text formatting
# so the literal decimal strings are preserved verbatim instead of being
# coerced into floating point numbers (matches the original inlineStr text cells).
$ws.Range("D2:E7").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = "Production_system"
$ws.Range("C2").Value = "svm"
$ws.Range("D2").Value = "0.171147107609372"
$ws.Range("E2").Value = "0.207770676738276"
$ws.Range("F2").Value = "0 12 18 66 21 3 1 1 0 11 51 57 4 0 1 0 0 5 16 37 33 8 2 0 0 6 26 66 24 4 1 0 0 2 19 59 19 0 0 0 0 0 8 33 144 65 1 0 0 2 16 66 20 7 2 0 1 2 12 28 57 16 1 0"

# Row 3
$ws.Range("B3").Value = "Freshness"
$ws.Range("C3").Value = "svm"
$ws.Range("D3").Value = "0.583380735800705"
$ws.Range("E3").Value = "0.717757875888754"
$ws.Range("F3").Value = "510 6 396 49"

# Row 4
$ws.Range("B4").Value = "Production_system"
$ws.Range("C4").Value = "svm"
$ws.Range("D4").Value = "0.173928420014363"
$ws.Range("E4").Value = "0.231955137665147"
$ws.Range("F4").Value = "0 0 4 25 9 0 0 1 0 1 32 62 18 0 1 0 0 0 5 43 44 2 1 0 0 0 6 27 7 0 0 0 0 0 1 55 39 0 0 0 0 0 5 24 145 61 1 0 0 1 2 33 4 0 0 0 1 0 1 41 56 13 1 0"

# Row 5
$ws.Range("B5").Value = "Freshness"
$ws.Range("C5").Value = "svm"
$ws.Range("D5").Value = "0.536465480987216"
$ws.Range("E5").Value = "0.678776482745529"
$ws.Range("F5").Value = "355 5 327 28"

# Row 6
$ws.Range("B6").Value = "Production_system"
$ws.Range("C6").Value = "svm"
$ws.Range("D6").Value = "0.223592118035333"
$ws.Range("E6").Value = "0.309105738378981"
$ws.Range("F6").Value = "2 5 41 47 22 3 0 1 24 69 11 6 1 0 1 0 0 2 11 31 47 5 0 0 0 0 6 45 33 33 0 0 0 9 5 35 53 3 0 0 0 18 23 49 101 46 0 0 0 1 6 45 25 44 0 0 1 6 3 28 66 6 0 0"

# Row 7
$ws.Range("B7").Value = "Freshness"
$ws.Range("C7").Value = "svm"
$ws.Range("D7").Value = "0.517770881278367"
$ws.Range("E7").Value = "0.667311578386005"
$ws.Range("F7").Value = "454 12 437 26"
